# Complete Georgia (GA) scraper - 25 staff members
$wb = $excel.ActiveWorkbook

# --- FL sheet: drop the leftover empty placeholder cells (D/F) ---
$flWs = $wb.Worksheets.Item("FL")
$flWs.Range("D4:D10").ClearContents()
$flWs.Range("F4").ClearContents()

# --- GA sheet: append the 25 scraped staff rows ---
$gaWs = $wb.Worksheets.Item("GA")
$gaWs.Range("A4").Value = "Tim Israel"
$gaWs.Range("B4").Value = "Director"
$gaWs.Range("A5").Value = "Cassia Baker"
$gaWs.Range("B5").Value = "Project Manager, Cybersecurity"
$gaWs.Range("A6").Value = "Michael Barker"
$gaWs.Range("B6").Value = "Project Manager, Cybersecurity"
$gaWs.Range("A7").Value = "Sam Darwin"
$gaWs.Range("B7").Value = "Project Manager, Process Improvement"
$gaWs.Range("A8").Value = "Sandra Enciso"
$gaWs.Range("B8").Value = "Project Manager, Energy and Sustainability"
$gaWs.Range("A9").Value = "Alfred Gardner"
$gaWs.Range("B9").Value = "Project Manager, Human Resources"
$gaWs.Range("A10").Value = "Bogna Grabicka"
$gaWs.Range("B10").Value = "Project Manager, Safety and Sustainability"
$gaWs.Range("A11").Value = "Kelly Grissom"
$gaWs.Range("B11").Value = "Project Manager, Energy and Sustainability"
$gaWs.Range("A12").Value = "Andy Helm"
$gaWs.Range("B12").Value = "Project Manager, Strategy and Leadership Development"
$gaWs.Range("A13").Value = "Dean Hettenbach"
$gaWs.Range("B13").Value = "Project Manager, Supply Chain and Technology"
$gaWs.Range("A14").Value = "Katie Hines"
$gaWs.Range("B14").Value = "Project Manager, Process Improvement"
$gaWs.Range("A15").Value = "Andrea Hines"
$gaWs.Range("B15").Value = "Project Manager, Food and Beverage"
$gaWs.Range("A16").Value = "Andrew Krejci"
$gaWs.Range("B16").Value = "Project Manager, Technology"
$gaWs.Range("A17").Value = "Ben Cheeks"
$gaWs.Range("B17").Value = "Region Manager, Coastal Georgia"
$gaWs.Range("A18").Value = "Jason Clarke"
$gaWs.Range("B18").Value = "Region Manager, Northeast Georgia"
$gaWs.Range("A19").Value = "Hank Hobbs"
$gaWs.Range("B19").Value = "Region Manager, South Georgia"
$gaWs.Range("A20").Value = "Paul LaVigna"
$gaWs.Range("B20").Value = "Region Manager, South Metro Atlanta"
$gaWs.Range("A21").Value = "Jay Boudreaux"
$gaWs.Range("B21").Value = "Senior Program and Operations Manager"
$gaWs.Range("A22").Value = "Anna Cali"
$gaWs.Range("B22").Value = "Instructional Systems Designer"
$gaWs.Range("A23").Value = "Jasmyn Green"
$gaWs.Range("B23").Value = "Program and Operations Manager"
$gaWs.Range("A24").Value = "Raine Hyde"
$gaWs.Range("B24").Value = "Marketing Strategist"
$gaWs.Range("A25").Value = "Megan Johnson"
$gaWs.Range("B25").Value = "Marketing Manager, Outreach"
$gaWs.Range("A26").Value = "Ieasha Jones"
$gaWs.Range("B26").Value = "Special Events"
$gaWs.Range("A27").Value = "Amber Kasselman"
$gaWs.Range("B27").Value = "Marketing Manager"
$gaWs.Range("A28").Value = "Caley Landau"
$gaWs.Range("B28").Value = "Marketing Strategist"
